$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '30.020.56'
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +7.62%  '

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.874.28'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +5.59%  '

# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.07%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '248.63'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +2.21%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.04%  '

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.4967'
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +1.62%  '

# Row 8
$ws.Cells.Item(8, 2).Value = 'OKB'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '45.73'
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +9.19%  '

# Row 9
$ws.Cells.Item(9, 2).Value = 'Cardano'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.2837'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +7.09%  '

# Row 10
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.06544'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +5.06%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'WrappedEther'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '1.871.19'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +5.26%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'Solana'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '17.07'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +4.96%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.07178'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +2.44%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.6597'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +7.12%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Litecoin'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '85.01'
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +7.14%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'Polkadot'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '4.793'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +4.19%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '30.007.84'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +7.61%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'Dai'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.9990'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.23%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'Avalanche'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '12.86'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +9.22%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '0.000007484'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +4.08%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'BinanceUSD'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.16%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '2.111.81'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +5.15%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '4.735'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +3.95%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'Cosmos'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '9.011'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +4.54%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'Chainlink'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '5.489'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +5.86%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '144.30'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +1.67%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'BitcoinCash'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '134.63'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +23.69%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '16.69'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +7.41%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '1.950'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +5.11%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '1.405'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.80%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '4.203'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.60%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '0.08583'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +4.32%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '3.877'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.98%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '0.05059'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +6.75%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '1.132'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +7.41%  '

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.03%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '0.6810'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +6.42%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'HuobiToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '2.701'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +4.13%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '2.322'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +13.77%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '2.737'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +5.92%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.9609'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.04%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.01624'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +5.92%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '6.053'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +3.12%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'PaxDollar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.01%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '102.90'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.65%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'TheSandbox'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.4176'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +6.26%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Aptos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '7.469'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +4.50%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Algorand'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.1250'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +5.03%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.05630'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +4.08%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Elrond'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '32.40'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +6.75%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '8.265'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +3.80%  '
